# "semana 52 de 2025" - add weekly columns 51 (BB) and 52 (BC) to the IRA
# extended weekly revision sheet, and refresh the two rows (38 and 57)
# whose weekly counts were recalculated/corrected for this release.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header labels for new weeks 51 and 52 (text, matching existing header style)
$c = $ws.Cells.Item(1, 54)
$c.NumberFormat = "@"
$c.Value = "51"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c = $ws.Cells.Item(1, 55)
$c.NumberFormat = "@"
$c.Value = "52"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108

# Row 2
$ws.Cells.Item(2, 54).Value = 32
$ws.Cells.Item(2, 55).Value = 32

# Row 3
$ws.Cells.Item(3, 54).Value = 41
$ws.Cells.Item(3, 55).Value = 32

# Row 5
$ws.Cells.Item(5, 54).Value = 3
$ws.Cells.Item(5, 55).Value = 1

# Row 6
$ws.Cells.Item(6, 54).Value = 74
$ws.Cells.Item(6, 55).Value = 61

# Row 7
$ws.Cells.Item(7, 54).Value = 24
$ws.Cells.Item(7, 55).Value = 28

# Row 8
$ws.Cells.Item(8, 54).Value = 17
$ws.Cells.Item(8, 55).Value = 14

# Row 9
$ws.Cells.Item(9, 54).Value = 6
$ws.Cells.Item(9, 55).Value = 5

# Row 10
$ws.Cells.Item(10, 54).Value = 3
$ws.Cells.Item(10, 55).Value = 1

# Row 11
$ws.Cells.Item(11, 54).Value = 2

# Row 14
$ws.Cells.Item(14, 54).Value = 1
$ws.Cells.Item(14, 55).Value = 1

# Row 15
$ws.Cells.Item(15, 54).Value = 1

# Row 16
$ws.Cells.Item(16, 54).Value = 1
$ws.Cells.Item(16, 55).Value = 2

# Row 17
$ws.Cells.Item(17, 54).Value = 1

# Row 23
$ws.Cells.Item(23, 54).Value = 2

# Row 25
$ws.Cells.Item(25, 54).Value = 26
$ws.Cells.Item(25, 55).Value = 27

# Row 28
$ws.Cells.Item(28, 54).Value = 195
$ws.Cells.Item(28, 55).Value = 135

# Row 29
$ws.Cells.Item(29, 54).Value = 0
$ws.Cells.Item(29, 55).Value = 0

# Row 30
$ws.Cells.Item(30, 54).Value = 42
$ws.Cells.Item(30, 55).Value = 36

# Row 31
$ws.Cells.Item(31, 54).Value = 0
$ws.Cells.Item(31, 55).Value = 1

# Row 35
$ws.Cells.Item(35, 54).Value = 41
$ws.Cells.Item(35, 55).Value = 49

# Row 36
$ws.Cells.Item(36, 54).Value = 4
$ws.Cells.Item(36, 55).Value = 4

# Row 37
$ws.Cells.Item(37, 54).Value = 14
$ws.Cells.Item(37, 55).Value = 9

# Row 38
$ws.Cells.Item(38, 4).Value = 77
$ws.Cells.Item(38, 5).Value = 85
$ws.Cells.Item(38, 6).Value = 86
$ws.Cells.Item(38, 7).Value = 79
$ws.Cells.Item(38, 8).Value = 56
$ws.Cells.Item(38, 9).Value = 30
$ws.Cells.Item(38, 10).Value = 66
$ws.Cells.Item(38, 11).Value = 79
$ws.Cells.Item(38, 12).Value = 64
$ws.Cells.Item(38, 13).Value = 70
$ws.Cells.Item(38, 14).Value = 74
$ws.Cells.Item(38, 15).Value = 83
$ws.Cells.Item(38, 16).Value = 70
$ws.Cells.Item(38, 17).Value = 69
$ws.Cells.Item(38, 18).Value = 73
$ws.Cells.Item(38, 19).Value = 55
$ws.Cells.Item(38, 20).Value = 58
$ws.Cells.Item(38, 21).Value = 54
$ws.Cells.Item(38, 22).Value = 41
$ws.Cells.Item(38, 23).Value = 68
$ws.Cells.Item(38, 24).Value = 56
$ws.Cells.Item(38, 25).Value = 63
$ws.Cells.Item(38, 26).Value = 58
$ws.Cells.Item(38, 27).Value = 58
$ws.Cells.Item(38, 28).Value = 56
$ws.Cells.Item(38, 29).Value = 54
$ws.Cells.Item(38, 30).Value = 53
$ws.Cells.Item(38, 31).Value = 58
$ws.Cells.Item(38, 32).Value = 70
$ws.Cells.Item(38, 33).Value = 63
$ws.Cells.Item(38, 34).Value = 51
$ws.Cells.Item(38, 35).Value = 51
$ws.Cells.Item(38, 36).Value = 68
$ws.Cells.Item(38, 37).Value = 76
$ws.Cells.Item(38, 38).Value = 57
$ws.Cells.Item(38, 39).Value = 57
$ws.Cells.Item(38, 40).Value = 47
$ws.Cells.Item(38, 41).Value = 55
$ws.Cells.Item(38, 42).Value = 48
$ws.Cells.Item(38, 43).Value = 51
$ws.Cells.Item(38, 44).Value = 55
$ws.Cells.Item(38, 45).Value = 54
$ws.Cells.Item(38, 46).Value = 42
$ws.Cells.Item(38, 47).Value = 59
$ws.Cells.Item(38, 48).Value = 50
$ws.Cells.Item(38, 49).Value = 44
$ws.Cells.Item(38, 50).Value = 59
$ws.Cells.Item(38, 51).Value = 64
$ws.Cells.Item(38, 52).Value = 55
$ws.Cells.Item(38, 53).Value = 51
$ws.Cells.Item(38, 54).Value = 69
$ws.Cells.Item(38, 55).Value = 83

# Row 41
$ws.Cells.Item(41, 54).Value = 8
$ws.Cells.Item(41, 55).Value = 6

# Row 42
$ws.Cells.Item(42, 54).Value = 115
$ws.Cells.Item(42, 55).Value = 51

# Row 43
$ws.Cells.Item(43, 54).Value = 64
$ws.Cells.Item(43, 55).Value = 39

# Row 45
$ws.Cells.Item(45, 54).Value = 24
$ws.Cells.Item(45, 55).Value = 269

# Row 46
$ws.Cells.Item(46, 54).Value = 95
$ws.Cells.Item(46, 55).Value = 63

# Row 47
$ws.Cells.Item(47, 54).Value = 175
$ws.Cells.Item(47, 55).Value = 111

# Row 48
$ws.Cells.Item(48, 54).Value = 7
$ws.Cells.Item(48, 55).Value = 6

# Row 49
$ws.Cells.Item(49, 54).Value = 120
$ws.Cells.Item(49, 55).Value = 85

# Row 50
$ws.Cells.Item(50, 54).Value = 1
$ws.Cells.Item(50, 55).Value = 2

# Row 51
$ws.Cells.Item(51, 54).Value = 0
$ws.Cells.Item(51, 55).Value = 0

# Row 54
$ws.Cells.Item(54, 54).Value = 13
$ws.Cells.Item(54, 55).Value = 7

# Row 55
$ws.Cells.Item(55, 54).Value = 1
$ws.Cells.Item(55, 55).Value = 4

# Row 56
$ws.Cells.Item(56, 54).Value = 1
$ws.Cells.Item(56, 55).Value = 0

# Row 57
$ws.Cells.Item(57, 15).Value = 11
$ws.Cells.Item(57, 17).Value = 8
$ws.Cells.Item(57, 18).Value = 16
$ws.Cells.Item(57, 20).Value = 7
$ws.Cells.Item(57, 22).Value = 6
$ws.Cells.Item(57, 24).Value = 7
$ws.Cells.Item(57, 25).Value = 3
$ws.Cells.Item(57, 29).Value = 4
$ws.Cells.Item(57, 36).Value = 12
$ws.Cells.Item(57, 40).Value = 4
$ws.Cells.Item(57, 41).Value = 5
$ws.Cells.Item(57, 42).Value = 6
$ws.Cells.Item(57, 54).Value = 4
$ws.Cells.Item(57, 55).Value = 6

# Row 58
$ws.Cells.Item(58, 54).Value = 16
$ws.Cells.Item(58, 55).Value = 16

# Row 59
$ws.Cells.Item(59, 54).Value = 9
$ws.Cells.Item(59, 55).Value = 10
